$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $origStyle = $ws.Range($addr).Style
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = $origStyle
}

# Row 2
Set-TextValue $ws "D2" "63.410.66"
$ws.Range("E2").Value = "  -1.04%  "

# Row 3
Set-TextValue $ws "D3" "2.718.22"
$ws.Range("E3").Value = "  -1.44%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
Set-TextValue $ws "D5" "561.01"
$ws.Range("E5").Value = "  -2.54%  "

# Row 6
Set-TextValue $ws "D6" "156.98"
$ws.Range("E6").Value = "  -1.20%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
Set-TextValue $ws "D8" "0.592"
$ws.Range("E8").Value = "  -2.03%  "

# Row 9
$ws.Range("E9").Value = "  -2.61%  "

# Row 10
Set-TextValue $ws "D10" "0.165"
$ws.Range("E10").Value = "  -0.58%  "

# Row 11
Set-TextValue $ws "D11" "5.59"
$ws.Range("E11").Value = "  -1.90%  "

# Row 12
Set-TextValue $ws "D12" "0.372"
$ws.Range("E12").Value = "  -3.98%  "

# Row 13
Set-TextValue $ws "D13" "3.195.55"
$ws.Range("E13").Value = "  -1.60%  "

# Row 14
Set-TextValue $ws "D14" "26.49"
$ws.Range("E14").Value = "  -1.63%  "

# Row 15
Set-TextValue $ws "D15" "63.287.91"
$ws.Range("E15").Value = "  -0.62%  "

# Row 16
$ws.Range("E16").Value = "  -3.01%  "

# Row 17
Set-TextValue $ws "D17" "2.716.54"
$ws.Range("E17").Value = "  -1.66%  "

# Row 18
Set-TextValue $ws "D18" "12.20"
$ws.Range("E18").Value = "  +0.61%  "

# Row 19
Set-TextValue $ws "D19" "4.68"
$ws.Range("E19").Value = "  -3.92%  "

# Row 20
Set-TextValue $ws "D20" "351.55"
$ws.Range("E20").Value = "  -1.68%  "

# Row 21
Set-TextValue $ws "D21" "6.48"
$ws.Range("E21").Value = "  -3.94%  "

# Row 22
Set-TextValue $ws "D22" "0.999"
$ws.Range("E22").Value = "  -0.01%  "

# Row 23
Set-TextValue $ws "D23" "0.513"
$ws.Range("E23").Value = "  -4.08%  "

# Row 24
Set-TextValue $ws "D24" "64.40"
$ws.Range("E24").Value = "  -1.79%  "

# Row 25
$ws.Range("E25").Value = "  -0.99%  "

# Row 26
$ws.Range("E26").Value = "  +0.02%  "

# Row 27
Set-TextValue $ws "D27" "8.21"
$ws.Range("E27").Value = "  -4.43%  "

# Row 28
Set-TextValue $ws "D28" "0.0₃0893"
$ws.Range("E28").Value = "  -1.88%  "

# Row 29
Set-TextValue $ws "D29" "1.38"
$ws.Range("E29").Value = "  +10.75%  "

# Row 30
$ws.Range("E30").Value = "  -0.69%  "

# Row 31
Set-TextValue $ws "D31" "7.19"
$ws.Range("E31").Value = "  -1.46%  "

# Row 32
Set-TextValue $ws "D32" "166.21"
$ws.Range("E32").Value = "  -1.75%  "

# Row 33
$ws.Range("E33").Value = "  -0.40%  "

# Row 34
Set-TextValue $ws "B34" "EthereumClassic"
Set-TextValue $ws "C34" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D34" "19.88"
$ws.Range("E34").Value = "  -1.91%  "

# Row 35
Set-TextValue $ws "B35" "USDe"
Set-TextValue $ws "C35" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws "D35" "0.999"
$ws.Range("E35").Value = "  +0.03%  "

# Row 36
Set-TextValue $ws "D36" "4.84"
$ws.Range("E36").Value = "  -1.74%  "

# Row 37
Set-TextValue $ws "D37" "1.78"
$ws.Range("E37").Value = "  -2.00%  "

# Row 38
Set-TextValue $ws "D38" "346.47"
$ws.Range("E38").Value = "  +0.13%  "

# Row 39
Set-TextValue $ws "D39" "0.961"
$ws.Range("E39").Value = "  -4.41%  "

# Row 40
Set-TextValue $ws "D40" "6.12"
$ws.Range("E40").Value = "  -3.30%  "

# Row 41
Set-TextValue $ws "D41" "4.06"
$ws.Range("E41").Value = "  -3.47%  "

# Row 42
Set-TextValue $ws "D42" "38.45"
$ws.Range("E42").Value = "  -1.91%  "

# Row 43
Set-TextValue $ws "D43" "21.45"
$ws.Range("E43").Value = "  -1.81%  "

# Row 44
Set-TextValue $ws "D44" "20.81"
$ws.Range("E44").Value = "  -2.95%  "

# Row 45
Set-TextValue $ws "D45" "0.0575"
$ws.Range("E45").Value = "  -2.71%  "

# Row 46
Set-TextValue $ws "D46" "0.625"
$ws.Range("E46").Value = "  -1.11%  "

# Row 47
Set-TextValue $ws "B47" "Aave"
Set-TextValue $ws "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D47" "131.97"
$ws.Range("E47").Value = "  -2.88%  "

# Row 48
Set-TextValue $ws "B48" "FirstDigitalUSD"
Set-TextValue $ws "C48" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D48" "0.998"
$ws.Range("E48").Value = "  -0.03%  "

# Row 49
$ws.Range("E49").Value = "  -3.51%  "

# Row 50
Set-TextValue $ws "B50" "Stellar"
Set-TextValue $ws "C50" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws "D50" "0.0985"
$ws.Range("E50").Value = "  -3.50%  "

# Row 51
Set-TextValue $ws "B51" "WhiteBITCoin"
Set-TextValue $ws "C51" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws "D51" "11.05"
$ws.Range("E51").Value = "  +0.14%  "
